$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 1871.1666
$ws.Range("I94").Value = 1871.1666
$ws.Range("K94").Value = 1871.1666
$ws.Range("M94").Value = -1420.1666

$ws.Range("H137").Value = 4008526
$ws.Range("I137").Value = 6420.8823
$ws.Range("K137").Value = 19262.6469
$ws.Range("M137").Value = -16712.6469

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 762.44446
$ws.Range("I4").Value = 1081
$ws.Range("J4").Value = 507.6
$ws.Range("K4").Value = 1081
$ws.Range("L4").Value = 507.6
$ws.Range("M4").Value = -965
$ws.Range("N4").Value = -739.6

$ws.Range("H32").Value = 215233.72
$ws.Range("I32").Value = 264236.25
$ws.Range("K32").Value = 264236.25
$ws.Range("M32").Value = -263949.25

$ws.Range("H37").Value = 49248.75
$ws.Range("I37").Value = 44998.5
$ws.Range("K37").Value = 44998.5
$ws.Range("M37").Value = -44725.5

$ws.Range("H102").Value = 2071.1904
$ws.Range("I102").Value = 2122.4
$ws.Range("K102").Value = 2122.4
$ws.Range("M102").Value = -500.4000000000001

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H132").Value = 2665.5217
$ws.Range("I132").Value = 2637.1875
$ws.Range("K132").Value = 7911.5625
$ws.Range("M132").Value = -5381.5625

$ws.Range("H134").Value = 139990
$ws.Range("J134").Value = 139990
$ws.Range("L134").Value = 139990
$ws.Range("N134").Value = -150130

$ws.Range("H137").Value = 105993.336
$ws.Range("J137").Value = 105993.336
$ws.Range("L137").Value = 105993.336
$ws.Range("N137").Value = -116193.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4398.727
$ws.Range("I94").Value = 2704.5715
$ws.Range("K94").Value = 2704.5715
$ws.Range("M94").Value = -2253.5715

$ws.Range("H99").Value = 14041.7
$ws.Range("I99").Value = 21666.334
$ws.Range("K99").Value = 21666.334
$ws.Range("M99").Value = -20168.334

$ws.Range("H134").Value = 30002212
$ws.Range("I134").Value = 1764.2084
$ws.Range("K134").Value = 5292.6252
$ws.Range("M134").Value = -2757.6252

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2424.5
$ws.Range("I31").Value = 2723.9092
$ws.Range("K31").Value = 2723.9092
$ws.Range("M31").Value = -2428.9092

$ws.Range("H34").Value = 2424.5
$ws.Range("I34").Value = 2723.9092
$ws.Range("K34").Value = 2723.9092
$ws.Range("M34").Value = -2521.9092

$ws.Range("H62").Value = 3533.762
$ws.Range("J62").Value = 3495.2
$ws.Range("L62").Value = 3495.2
$ws.Range("N62").Value = -4743.2

$ws.Range("H65").Value = 3533.762
$ws.Range("J65").Value = 3495.2
$ws.Range("L65").Value = 17476
$ws.Range("N65").Value = -23716

$ws.Range("H99").Value = 35716892
$ws.Range("I99").Value = 2629.1
$ws.Range("K99").Value = 2629.1
$ws.Range("M99").Value = -1131.1

$ws.Range("H122").Value = 3532.25
$ws.Range("I122").Value = 3186.3
$ws.Range("K122").Value = 9558.900000000001
$ws.Range("M122").Value = -7108.900000000001

$ws.Range("H126").Value = 35716892
$ws.Range("I126").Value = 2629.1
$ws.Range("K126").Value = 7887.299999999999
$ws.Range("M126").Value = -5417.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 14375.875
$ws.Range("I3").Value = 11502.833
$ws.Range("K3").Value = 34508.499
$ws.Range("M3").Value = -34396.499

$ws.Range("H18").Value = 957.4286
$ws.Range("I18").Value = 283.66666
$ws.Range("K18").Value = 850.9999799999999
$ws.Range("M18").Value = -681.9999799999999

$ws.Range("H26").Value = 194.75
$ws.Range("J26").Value = 190
$ws.Range("L26").Value = 570
$ws.Range("N26").Value = -1146

$ws.Range("H37").Value = 124992.5
$ws.Range("J37").Value = 124992.5
$ws.Range("L37").Value = 374977.5
$ws.Range("N37").Value = -375201.5

$ws.Range("H41").Value = 23200.727
$ws.Range("I41").Value = 900
$ws.Range("J41").Value = 35944
$ws.Range("K41").Value = 2700
$ws.Range("L41").Value = 107832
$ws.Range("M41").Value = -2362
$ws.Range("N41").Value = -108508

$ws.Range("H52").Value = 1636.5
$ws.Range("J52").Value = 1636.5
$ws.Range("L52").Value = 4909.5
$ws.Range("N52").Value = -5441.5

$ws.Range("H55").Value = 4849.625
$ws.Range("J55").Value = 5411.154
$ws.Range("L55").Value = 16233.462
$ws.Range("N55").Value = -16587.462

$ws.Range("H60").Value = 5000
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()

$ws.Range("H81").Value = 7358.7144
$ws.Range("J81").Value = 8302.200000000001
$ws.Range("L81").Value = 24906.6
$ws.Range("N81").Value = -27152.6

$ws.Range("H84").Value = 7358.7144
$ws.Range("J84").Value = 8302.200000000001
$ws.Range("L84").Value = 74719.8
$ws.Range("N84").Value = -85951.8

$ws.Range("H114").Value = 1400.8
$ws.Range("I114").Value = 1112.1666
$ws.Range("K114").Value = 3336.4998
$ws.Range("M114").Value = -82.49980000000005

$ws.Range("H136").Value = 9331
$ws.Range("I136").Value = 7331.3
$ws.Range("K136").Value = 21993.9
$ws.Range("M136").Value = -16893.9

$ws.Range("H137").Value = 2584.35
$ws.Range("I137").Value = 1543.6364
$ws.Range("K137").Value = 4630.9092
$ws.Range("M137").Value = 469.0907999999999

$ws.Range("H139").Value = 3949.2
$ws.Range("I139").Value = 2689.25
$ws.Range("K139").Value = 8067.75
$ws.Range("M139").Value = -2927.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4100.143
$ws.Range("I70").Value = 3854
$ws.Range("J70").Value = 5002.6665
$ws.Range("K70").Value = 3854
$ws.Range("L70").Value = 5002.6665
$ws.Range("M70").Value = -3584
$ws.Range("N70").Value = -5542.6665

$ws.Range("H73").Value = 4100.143
$ws.Range("I73").Value = 3854
$ws.Range("J73").Value = 5002.6665
$ws.Range("K73").Value = 3854
$ws.Range("L73").Value = 5002.6665
$ws.Range("M73").Value = -2918
$ws.Range("N73").Value = -6874.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H110").Value = 30500
$ws.Range("J110").Value = 30500
$ws.Range("L110").Value = 30500
$ws.Range("N110").Value = -38680

$ws.Range("H122").Value = 3550.0881
$ws.Range("I122").Value = 3082
$ws.Range("K122").Value = 9246
$ws.Range("M122").Value = -6796

$ws.Range("H132").Value = 2959.1853
$ws.Range("I132").Value = 2665.2354
$ws.Range("J132").Value = 3458.9
$ws.Range("K132").Value = 7995.706200000001
$ws.Range("L132").Value = 10376.7
$ws.Range("M132").Value = -5465.706200000001
$ws.Range("N132").Value = -15436.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 34815
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 34815
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 34815
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -35369

$ws.Range("H123").Value = 93747
$ws.Range("J123").Value = 93747
$ws.Range("L123").Value = 93747
$ws.Range("N123").Value = -103547

$ws.Range("H125").Value = 63711.145
$ws.Range("J125").Value = 63711.145
$ws.Range("L125").Value = 63711.145
$ws.Range("N125").Value = -73551.14499999999

$ws.Range("H132").Value = 2589.625
$ws.Range("I132").Value = 1871.4706
$ws.Range("K132").Value = 5614.4118
$ws.Range("M132").Value = -3084.4118

$ws.Range("H135").Value = 97560.664
$ws.Range("J135").Value = 97560.664
$ws.Range("L135").Value = 97560.664
$ws.Range("N135").Value = -107700.664
